# "Import to multiple tables and course addition"
#
# The sample class-list header used the column name "student_name"; the
# importer now expects a generic "name" column (so the same sheet shape can
# feed multiple destination tables). Rename the A1 header accordingly - the
# rest of the data (matric_no header and the two student rows) is untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "name"
